$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "LIDIA KATHERINE CHAVEZ REYES"
$ws.Range("B2").Value = 112

$ws.Range("A3").Value = "SANCHEZ LUCUMI DIANA"
$ws.Range("B3").Value = 112

$ws.Range("A4").Value = "TORRES ASCORBE CESAR RAMNCES"
$ws.Range("B4").Value = 109

$ws.Range("A5").Value = "SANCHEZ SALDAÑA FRANK REGINALDO"
$ws.Range("B5").Value = 97

$ws.Range("A6").Value = "BLANCO LOZANO ANDREA MILAGROS"
$ws.Range("B6").Value = 85

$ws.Range("A7").Value = "BURGA MEDINA SHIRLEY ROCIO"
$ws.Range("B7").Value = 85

$ws.Range("A8").Value = "DE LA CRUZ CARDENAS RUTH LUCERO"
$ws.Range("B8").Value = 84

$ws.Range("A9").Value = "URBINA ANTICONA ALEX BRUNO"
$ws.Range("B9").Value = 83
